$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Col1a2"
$ws.Range("C2").Value = "Itga2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 19.36022366666667
$ws.Range("H2").Value = 58.080671
$ws.Range("I2").Value = 0.005884129141485179
$ws.Range("J2").Value = 0.005884129141485179
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.740822
$ws.Range("N2").Value = 5.222466
$ws.Range("O2").Value = 0.4863878955914668
$ws.Range("P2").Value = 0.4863878955914669
$ws.Range("Q2").Value = 33.702703283854
$ws.Range("R2").Value = 303.324329554686
$ws.Range("S2").Value = 0.002861969190515401
$ws.Range("T2").Value = 0.002861969190515401

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Col1a2"
$ws.Range("C3").Value = "Itga2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 19.36022366666667
$ws.Range("H3").Value = 58.080671
$ws.Range("I3").Value = 0.005884129141485179
$ws.Range("J3").Value = 0.005884129141485179
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.077748
$ws.Range("N3").Value = 3.233244
$ws.Range("O3").Value = 0.3011241710513264
$ws.Range("P3").Value = 0.3011241710513265
$ws.Range("Q3").Value = 20.86544233630266
$ws.Range("R3").Value = 187.788981026724
$ws.Range("S3").Value = 0.001771853510088678
$ws.Range("T3").Value = 0.001771853510088678

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Col1a2"
$ws.Range("C4").Value = "Itga2"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 19.36022366666667
$ws.Range("H4").Value = 58.080671
$ws.Range("I4").Value = 0.005884129141485179
$ws.Range("J4").Value = 0.005884129141485179
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.03488166666666666
$ws.Range("N4").Value = 0.104645
$ws.Range("O4").Value = 0.009745982326006345
$ws.Range("P4").Value = 0.009745982326006345
$ws.Range("Q4").Value = 0.6753168685327777
$ws.Range("R4").Value = 6.077851816795
$ws.Range("S4").Value = 0.00005734661861685344
$ws.Range("T4").Value = 0.00005734661861685344

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Col1a2"
$ws.Range("C5").Value = "Itga2"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 19.36022366666667
$ws.Range("H5").Value = 58.080671
$ws.Range("I5").Value = 0.005884129141485179
$ws.Range("J5").Value = 0.005884129141485179
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.7256300000000001
$ws.Range("N5").Value = 2.17689
$ws.Range("O5").Value = 0.2027419510312003
$ws.Range("P5").Value = 0.2027419510312003
$ws.Range("Q5").Value = 14.04835909924333
$ws.Range("R5").Value = 126.43523189319
$ws.Range("S5").Value = 0.001192959822264247
$ws.Range("T5").Value = 0.001192959822264247

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Col1a2"
$ws.Range("C6").Value = "Itga2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3161.845459
$ws.Range("H6").Value = 9485.536377
$ws.Range("I6").Value = 0.9609758299542277
$ws.Range("J6").Value = 0.9609758299542278
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.740822
$ws.Range("N6").Value = 5.222466
$ws.Range("O6").Value = 0.4863878955914668
$ws.Range("P6").Value = 0.4863878955914669
$ws.Range("Q6").Value = 5504.210135627298
$ws.Range("R6").Value = 49537.89122064568
$ws.Range("S6").Value = 0.4674070116457001
$ws.Range("T6").Value = 0.4674070116457002

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Col1a2"
$ws.Range("C7").Value = "Itga2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3161.845459
$ws.Range("H7").Value = 9485.536377
$ws.Range("I7").Value = 0.9609758299542277
$ws.Range("J7").Value = 0.9609758299542278
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.077748
$ws.Range("N7").Value = 3.233244
$ws.Range("O7").Value = 0.3011241710513264
$ws.Range("P7").Value = 0.3011241710513265
$ws.Range("Q7").Value = 3407.672619746332
$ws.Range("R7").Value = 30669.05357771699
$ws.Range("S7").Value = 0.2893730501953273
$ws.Range("T7").Value = 0.2893730501953273

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Col1a2"
$ws.Range("C8").Value = "Itga2"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3161.845459
$ws.Range("H8").Value = 9485.536377
$ws.Range("I8").Value = 0.9609758299542277
$ws.Range("J8").Value = 0.9609758299542278
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.03488166666666666
$ws.Range("N8").Value = 0.104645
$ws.Range("O8").Value = 0.009745982326006345
$ws.Range("P8").Value = 0.009745982326006345
$ws.Range("Q8").Value = 110.2904393523517
$ws.Range("R8").Value = 992.613954171165
$ws.Range("S8").Value = 0.009365653454453182
$ws.Range("T8").Value = 0.009365653454453183

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Col1a2"
$ws.Range("C9").Value = "Itga2"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 3161.845459
$ws.Range("H9").Value = 9485.536377
$ws.Range("I9").Value = 0.9609758299542277
$ws.Range("J9").Value = 0.9609758299542278
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.7256300000000001
$ws.Range("N9").Value = 2.17689
$ws.Range("O9").Value = 0.2027419510312003
$ws.Range("P9").Value = 0.2027419510312003
$ws.Range("Q9").Value = 2294.32992041417
$ws.Range("R9").Value = 20648.96928372753
$ws.Range("S9").Value = 0.1948301146587471
$ws.Range("T9").Value = 0.1948301146587471

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Col1a2"
$ws.Range("C10").Value = "Itga2"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.055785333333333
$ws.Range("H10").Value = 6.167356
$ws.Range("I10").Value = 0.0006248123263850286
$ws.Range("J10").Value = 0.0006248123263850286
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.740822
$ws.Range("N10").Value = 5.222466
$ws.Range("O10").Value = 0.4863878955914668
$ws.Range("P10").Value = 0.4863878955914669
$ws.Range("Q10").Value = 3.578756335544
$ws.Range("R10").Value = 32.208807019896
$ws.Range("S10").Value = 0.0003039011525700228
$ws.Range("T10").Value = 0.0003039011525700228

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Col1a2"
$ws.Range("C11").Value = "Itga2"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 2.055785333333333
$ws.Range("H11").Value = 6.167356
$ws.Range("I11").Value = 0.0006248123263850286
$ws.Range("J11").Value = 0.0006248123263850286
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 1.077748
$ws.Range("N11").Value = 3.233244
$ws.Range("O11").Value = 0.3011241710513264
$ws.Range("P11").Value = 0.3011241710513265
$ws.Range("Q11").Value = 2.215618531429333
$ws.Range("R11").Value = 19.940566782864
$ws.Range("S11").Value = 0.0001881460938453426
$ws.Range("T11").Value = 0.0001881460938453426

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Col1a2"
$ws.Range("C12").Value = "Itga2"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 2.055785333333333
$ws.Range("H12").Value = 6.167356
$ws.Range("I12").Value = 0.0006248123263850286
$ws.Range("J12").Value = 0.0006248123263850286
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.03488166666666666
$ws.Range("N12").Value = 0.104645
$ws.Range("O12").Value = 0.009745982326006345
$ws.Range("P12").Value = 0.009745982326006345
$ws.Range("Q12").Value = 0.07170921873555554
$ws.Range("R12").Value = 0.64538296862
$ws.Range("S12").Value = 0.000006089409890019397
$ws.Range("T12").Value = 0.000006089409890019397

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Col1a2"
$ws.Range("C13").Value = "Itga2"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 2.055785333333333
$ws.Range("H13").Value = 6.167356
$ws.Range("I13").Value = 0.0006248123263850286
$ws.Range("J13").Value = 0.0006248123263850286
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.7256300000000001
$ws.Range("N13").Value = 2.17689
$ws.Range("O13").Value = 0.2027419510312003
$ws.Range("P13").Value = 0.2027419510312003
$ws.Range("Q13").Value = 1.491739511426667
$ws.Range("R13").Value = 13.42565560284
$ws.Range("S13").Value = 0.0001266756700796438
$ws.Range("T13").Value = 0.0001266756700796438

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Col1a2"
$ws.Range("C14").Value = "Itga2"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 106.9830526666667
$ws.Range("H14").Value = 320.949158
$ws.Range("I14").Value = 0.03251522857790212
$ws.Range("J14").Value = 0.03251522857790212
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 1.740822
$ws.Range("N14").Value = 5.222466
$ws.Range("O14").Value = 0.4863878955914668
$ws.Range("P14").Value = 0.4863878955914669
$ws.Range("Q14").Value = 186.238451709292
$ws.Range("R14").Value = 1676.146065383628
$ws.Range("S14").Value = 0.01581501360268133
$ws.Range("T14").Value = 0.01581501360268134

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Col1a2"
$ws.Range("C15").Value = "Itga2"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 106.9830526666667
$ws.Range("H15").Value = 320.949158
$ws.Range("I15").Value = 0.03251522857790212
$ws.Range("J15").Value = 0.03251522857790212
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 1.077748
$ws.Range("N15").Value = 3.233244
$ws.Range("O15").Value = 0.3011241710513264
$ws.Range("P15").Value = 0.3011241710513265
$ws.Range("Q15").Value = 115.3007710453947
$ws.Range("R15").Value = 1037.706939408552
$ws.Range("S15").Value = 0.009791121252065176
$ws.Range("T15").Value = 0.009791121252065178

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Col1a2"
$ws.Range("C16").Value = "Itga2"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 106.9830526666667
$ws.Range("H16").Value = 320.949158
$ws.Range("I16").Value = 0.03251522857790212
$ws.Range("J16").Value = 0.03251522857790212
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.03488166666666666
$ws.Range("N16").Value = 0.104645
$ws.Range("O16").Value = 0.009745982326006345
$ws.Range("P16").Value = 0.009745982326006345
$ws.Range("Q16").Value = 3.731747182101111
$ws.Range("R16").Value = 33.58572463891
$ws.Range("S16").Value = 0.0003168928430462905
$ws.Range("T16").Value = 0.0003168928430462905

# Row 17
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Col1a2"
$ws.Range("C17").Value = "Itga2"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 106.9830526666667
$ws.Range("H17").Value = 320.949158
$ws.Range("I17").Value = 0.03251522857790212
$ws.Range("J17").Value = 0.03251522857790212
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.7256300000000001
$ws.Range("N17").Value = 2.17689
$ws.Range("O17").Value = 0.2027419510312003
$ws.Range("P17").Value = 0.2027419510312003
$ws.Range("Q17").Value = 77.63011250651334
$ws.Range("R17").Value = 698.6710125586201
$ws.Range("S17").Value = 0.006592200880109316
$ws.Range("T17").Value = 0.006592200880109316
